# Add a new rule row (CONVERT_WEIGHT for 홍게) to the OptionRules sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OptionRules")

# Switch to / activate the OptionRules tab, as the user was working there.
$ws.Activate()

# Insert a new row above row 12, shifting rule rows 12-33 down to 13-34.
$ws.Rows(12).Insert()

# Copy the formatting from the row above (row 11) into the newly inserted
# row so the new row matches the rest of the table's style.
$ws.Range("A11:F11").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Renumber the "순서" (sequence) column for every rule row that got shifted
# down by the insert (they were typed values, not a formula).
for ($r = 34; $r -ge 13; $r--) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r, 1).Value + 1
}

# Fill in the new rule row.
$ws.Cells.Item(12, 1).Value = 12
$ws.Cells.Item(12, 2).Value = "ALL"
$ws.Cells.Item(12, 3).Value = "홍게"
$ws.Cells.Item(12, 4).Value = "CONVERT_WEIGHT"
$ws.Cells.Item(12, 5).Value = "kg"
$ws.Cells.Item(12, 6).Value = "g을 kg로 변환하고 수량 곱함 (일반규칙)"

# Restore the view: scrolled down a bit, with A33 selected.
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("A33").Select()

Write-Host "done"
